$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.389.53"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.848.69"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6296"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07649"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.49"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07748"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "1.838.88"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.010"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001088"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6794"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.54"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "2.091.83"
$ws.Range("E17").Value = "  -7.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.137"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "29.417.43"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.446"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.18"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.383"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.66"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.471"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.303"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05637"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.045"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.158"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7096"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.774"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").Value = "1.228.67"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01799"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.454"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9134"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.0000"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "2.001.18"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.13"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.159"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4012"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.034"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.689"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.30%  "
